$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 110: 2025-10-25 (serial 45955), 四方坪站充电量(kw)
$ws.Cells.Item(110, 1).Value = 45955
$ws.Cells.Item(110, 2).Value = "四方坪站充电量(kw)"
$row110 = @(
    900.97700000000009, 1702.6160000000004, 374.84000000000003, 618.02599999999984,
    385.91400000000004, 655.79700000000014, 463.88899999999995, 247.92599999999996,
    80.251000000000005, 168.73000000000002, 222.25, 263.529, 782.66300000000012,
    1338.4579999999999, 747.38, 354.37400000000002, 239.41900000000004, 267.00300000000004,
    92.65, 173.49, 47.365000000000002, 62.66, 85.97, 133.22399999999999
)
for ($i = 0; $i -lt $row110.Length; $i++) {
    $ws.Cells.Item(110, 3 + $i).Value = $row110[$i]
}

# Row 111: 2025-10-25 (serial 45955), 高岭站充电量(kw)
$ws.Cells.Item(111, 1).Value = 45955
$ws.Cells.Item(111, 2).Value = "高岭站充电量(kw)"
$row111 = @(
    363.14699999999999, 331.78399999999999, 384.15, 44.338999999999999,
    67.266000000000005, 202.733, 374.48099999999999, 113.149,
    389.02600000000007, 92.569000000000003, 186.673, 352.68700000000001,
    416.73399999999998, 516.01299999999992, 258.74599999999998, 334.11500000000001,
    250.09199999999998, 182.23299999999998, 0, 27.478000000000002,
    138.542, 22.335999999999999, 54.697000000000003, 11.282
)
for ($i = 0; $i -lt $row111.Length; $i++) {
    $ws.Cells.Item(111, 3 + $i).Value = $row111[$i]
}

# Match the author's final selection from the edit
$ws.Range("J114").Select()
